$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Column widths
#    Column A -> as close as possible to 37.42578125 (grid snaps to /6, 37.5 is closest)
#    Column B -> 21 exactly
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 36.67
$ws.Columns.Item(2).ColumnWidth = 20.17

# ---------------------------------------------------------------------------
# 2. New shared strings are typed in this order in the source workbook:
#    B22, B24, A20, A21, K2
# ---------------------------------------------------------------------------

# --- Row 22/23/24 informational block -------------------------------------
$ws.Range("B22:I23").Merge()
$ws.Range("B22").Value = "Une intersection : l'impact dans la matrice de satisfaction globale"
$ws.Range("B22:I23").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B22:I23").VerticalAlignment = -4108     # xlCenter

$ws.Range("B24:I24").Merge()
$ws.Range("B24").Value = "CHACUN DES COLONNES DONNE UN SCORE DE 1"
$ws.Range("B24:I24").HorizontalAlignment = -4108
$ws.Range("B24:I24").VerticalAlignment = -4108

# --- Row 20: TOTAL (1) ------------------------------------------------------
$ws.Range("A20").Value = "TOTAL (1)"

$rng20 = $ws.Range("B20:J20")
$rng20.Value = 1
$rng20.Font.Bold = $true
$rng20.HorizontalAlignment = -4108
$rng20.VerticalAlignment = -4108

# --- Row 21: PONDERATION (variable contexte) --------------------------------
$ws.Rows.Item(21).RowHeight = 43.5

$ws.Range("A3").Copy($ws.Range("A21"))
$ws.Range("A21").Value = "PONDÉRATION (Variable contexte)"
$ws.Range("A21").VerticalAlignment = -4108

$rng21 = $ws.Range("B21:K21")
$rng21.Font.Bold = $true
$rng21.HorizontalAlignment = -4108
$rng21.VerticalAlignment = -4108

$ws.Range("C11").Copy($ws.Range("E21"))
$ws.Range("E21").Font.Bold = $true
$ws.Range("E21").HorizontalAlignment = -4108
$ws.Range("E21").VerticalAlignment = -4108

$ws.Range("K21").Value = 1
$ws.Range("K21").Font.Bold = $true
$ws.Range("K21").HorizontalAlignment = -4108
$ws.Range("K21").VerticalAlignment = -4108

# --- Row 2: K2 total label ---------------------------------------------------
$ws.Range("K2").Value = "TOTAL"

# ---------------------------------------------------------------------------
# 3. Existing value tweaks
# ---------------------------------------------------------------------------
$ws.Range("F12").Copy($ws.Range("C7"))     # "0.2" text, same style as before
$ws.Range("F6").Copy($ws.Range("D18"))     # "0.1" text
$ws.Range("C11").Copy($ws.Range("C19"))    # "0.3" text

# ---------------------------------------------------------------------------
# 4. Row 54 (extends the used range / dimension down to row 54)
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A54").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 5. Selection cosmetics
# ---------------------------------------------------------------------------
$ws.Range("J29").Select()
